$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing "source" block (MTICM / Ministry...) down to rows 25-26 ---
# Capture current text before it gets overwritten by the new table rows.
$mticm = $ws.Range("A19").Value()
$ministry = $ws.Range("A20").Value()

$ws.Range("A25").Value = $mticm
$ws.Range("A25").Style = "title"

$ws.Range("A26").Value = $ministry
$ws.Range("A26").Style = "source"

# --- New header row (row 16): Number of employees / Assets / Turnover ---
$ws.Range("B16").Value = "Number of employees"
$ws.Range("C16").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D16").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B16:D16").Style = "title"

# --- New data rows 17-20: size category / employee range / blank assets / blank turnover ---
$ws.Range("A17:D20").Style = "Normal"

$ws.Range("A17").Value = "Micro"
$ws.Range("B17").Value = "<3"

$ws.Range("A18").Value = "Small"
$ws.Range("B18").Value = "3-9"

$ws.Range("A19").Value = "Medium"
$ws.Range("B19").Value = "10-49"

$ws.Range("A20").Value = "Large"
$ws.Range("B20").Value = ">=50"
